$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.539.46"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "2.109.68"
$ws.Range("E3").Value = "  +0.91%  "
$r = $ws.Cells.Item(4, 4)
$r.NumberFormat = "@"
$r.Value = "1.010"
$r.Style = $ws.Cells.Item(4, 2).Style
$r = $ws.Cells.Item(5, 4)
$r.NumberFormat = "@"
$r.Value = "335.73"
$r.Style = $ws.Cells.Item(5, 2).Style
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("E6").Value = "  +0.63%  "
$r = $ws.Cells.Item(7, 4)
$r.NumberFormat = "@"
$r.Value = "0.5230"
$r.Style = $ws.Cells.Item(7, 2).Style
$ws.Range("E7").Value = "  +0.60%  "
$r = $ws.Cells.Item(8, 4)
$r.NumberFormat = "@"
$r.Value = "0.4543"
$r.Style = $ws.Cells.Item(8, 2).Style
$ws.Range("E8").Value = "  +4.20%  "
$r = $ws.Cells.Item(9, 4)
$r.NumberFormat = "@"
$r.Value = "55.24"
$r.Style = $ws.Cells.Item(9, 2).Style
$ws.Range("E9").Value = "  +2.71%  "
$r = $ws.Cells.Item(10, 4)
$r.NumberFormat = "@"
$r.Value = "0.09094"
$r.Style = $ws.Cells.Item(10, 2).Style
$ws.Range("E10").Value = "  +2.99%  "
$ws.Range("E11").Value = "  +1.55%  "
$r = $ws.Cells.Item(12, 4)
$r.NumberFormat = "@"
$r.Value = "24.52"
$r.Style = $ws.Cells.Item(12, 2).Style
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("D13").Value = "2.103.65"
$ws.Range("E13").Value = "  +0.91%  "
$r = $ws.Cells.Item(14, 4)
$r.NumberFormat = "@"
$r.Value = "6.826"
$r.Style = $ws.Cells.Item(14, 2).Style
$ws.Range("E14").Value = "  +2.33%  "
$r = $ws.Cells.Item(15, 4)
$r.NumberFormat = "@"
$r.Value = "8.100"
$r.Style = $ws.Cells.Item(15, 2).Style
$ws.Range("E15").Value = "  +5.81%  "
$r = $ws.Cells.Item(16, 4)
$r.NumberFormat = "@"
$r.Value = "0.00001173"
$r.Style = $ws.Cells.Item(16, 2).Style
$ws.Range("E16").Value = "  +4.91%  "
$r = $ws.Cells.Item(17, 4)
$r.NumberFormat = "@"
$r.Value = "96.95"
$r.Style = $ws.Cells.Item(17, 2).Style
$ws.Range("E17").Value = "  +1.41%  "
$ws.Range("E18").Value = "  +0.57%  "
$r = $ws.Cells.Item(19, 4)
$r.NumberFormat = "@"
$r.Value = "0.06677"
$r.Style = $ws.Cells.Item(19, 2).Style
$ws.Range("E19").Value = "  +1.37%  "
$r = $ws.Cells.Item(20, 4)
$r.NumberFormat = "@"
$r.Value = "19.35"
$r.Style = $ws.Cells.Item(20, 2).Style
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "30.595.44"
$ws.Range("E23").Value = "  +0.34%  "
$r = $ws.Cells.Item(24, 4)
$r.NumberFormat = "@"
$r.Value = "12.76"
$r.Style = $ws.Cells.Item(24, 2).Style
$ws.Range("E24").Value = "  +4.63%  "
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("D26").Value = "2.351.87"
$ws.Range("E26").Value = "  +0.84%  "
$r = $ws.Cells.Item(27, 4)
$r.NumberFormat = "@"
$r.Value = "22.24"
$r.Style = $ws.Cells.Item(27, 2).Style
$ws.Range("E27").Value = "  +0.31%  "
$r = $ws.Cells.Item(28, 4)
$r.NumberFormat = "@"
$r.Value = "163.64"
$r.Style = $ws.Cells.Item(28, 2).Style
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("E29").Value = "  -1.40%  "
$r = $ws.Cells.Item(30, 4)
$r.NumberFormat = "@"
$r.Value = "133.45"
$r.Style = $ws.Cells.Item(30, 2).Style
$ws.Range("E30").Value = "  +1.63%  "
$r = $ws.Cells.Item(31, 4)
$r.NumberFormat = "@"
$r.Value = "1.209"
$r.Style = $ws.Cells.Item(31, 2).Style
$ws.Range("E31").Value = "  +2.29%  "
$ws.Range("E32").Value = "  +0.17%  "
$r = $ws.Cells.Item(33, 4)
$r.NumberFormat = "@"
$r.Value = "1.635"
$r.Style = $ws.Cells.Item(33, 2).Style
$ws.Range("E33").Value = "  -0.59%  "
$r = $ws.Cells.Item(34, 4)
$r.NumberFormat = "@"
$r.Value = "6.348"
$r.Style = $ws.Cells.Item(34, 2).Style
$ws.Range("E34").Value = "  +3.50%  "
$r = $ws.Cells.Item(35, 4)
$r.NumberFormat = "@"
$r.Value = "3.957"
$r.Style = $ws.Cells.Item(35, 2).Style
$ws.Range("E35").Value = "  +1.32%  "
$r = $ws.Cells.Item(36, 4)
$r.NumberFormat = "@"
$r.Value = "10.43"
$r.Style = $ws.Cells.Item(36, 2).Style
$ws.Range("E36").Value = "  +3.14%  "
$r = $ws.Cells.Item(37, 4)
$r.NumberFormat = "@"
$r.Value = "5.901"
$r.Style = $ws.Cells.Item(37, 2).Style
$r = $ws.Cells.Item(38, 4)
$r.NumberFormat = "@"
$r.Value = "0.02611"
$r.Style = $ws.Cells.Item(38, 2).Style
$ws.Range("E38").Value = "  +1.56%  "
$r = $ws.Cells.Item(39, 4)
$r.NumberFormat = "@"
$r.Value = "0.06800"
$r.Style = $ws.Cells.Item(39, 2).Style
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("E40").Value = "  +2.99%  "
$r = $ws.Cells.Item(41, 4)
$r.NumberFormat = "@"
$r.Value = "12.56"
$r.Style = $ws.Cells.Item(41, 2).Style
$ws.Range("E41").Value = "  -0.21%  "
$r = $ws.Cells.Item(42, 4)
$r.NumberFormat = "@"
$r.Value = "0.6848"
$r.Style = $ws.Cells.Item(42, 2).Style
$ws.Range("E42").Value = "  -0.16%  "
$r = $ws.Cells.Item(43, 4)
$r.NumberFormat = "@"
$r.Value = "1.259"
$r.Style = $ws.Cells.Item(43, 2).Style
$ws.Range("E43").Value = "  +0.16%  "
$r = $ws.Cells.Item(44, 4)
$r.NumberFormat = "@"
$r.Value = "14.59"
$r.Style = $ws.Cells.Item(44, 2).Style
$ws.Range("E44").Value = "  +5.51%  "
$r = $ws.Cells.Item(45, 4)
$r.NumberFormat = "@"
$r.Value = "0.6426"
$r.Style = $ws.Cells.Item(45, 2).Style
$ws.Range("E45").Value = "  +1.67%  "
$r = $ws.Cells.Item(46, 4)
$r.NumberFormat = "@"
$r.Value = "2.293"
$r.Style = $ws.Cells.Item(46, 2).Style
$ws.Range("E46").Value = "  +4.80%  "
$r = $ws.Cells.Item(47, 4)
$r.NumberFormat = "@"
$r.Value = "3.683"
$r.Style = $ws.Cells.Item(47, 2).Style
$ws.Range("E47").Value = "  +1.57%  "
$r = $ws.Cells.Item(48, 4)
$r.NumberFormat = "@"
$r.Value = "0.00000000352"
$r.Style = $ws.Cells.Item(48, 2).Style
$ws.Range("E48").Value = "  +18.76%  "
$ws.Range("E49").Value = "  +1.03%  "
$r = $ws.Cells.Item(50, 4)
$r.NumberFormat = "@"
$r.Value = "83.00"
$r.Style = $ws.Cells.Item(50, 2).Style
$ws.Range("E50").Value = "  +1.83%  "
$r = $ws.Cells.Item(51, 4)
$r.NumberFormat = "@"
$r.Value = "0.3368"
$r.Style = $ws.Cells.Item(51, 2).Style
$ws.Range("E51").Value = "  +12.72%  "